$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) Append a new sentence to the end of the paragraph about the
#    three-toed woodpecker being dependent on continuous old-growth
#    forest.
# ---------------------------------------------------------------------
$oldTail = "Virkkala, 1991)."
$newTail = "Virkkala, 1991). Det finns studier som visar att naturskogens självgallringsprocess som den tretåiga hackspetten är beroende av, inte kan ersättas med efterlämnad hänsyn i skogsbruket (Imbeau & Desrochers, 2002)."
$d.Content.Find.Execute($oldTail, $true, $false, $false, $false, $false, $true, 1, $false, $newTail, 2) | Out-Null

# ---------------------------------------------------------------------
# 2) Insert a new reference paragraph (Imbeau & Desrochers, 2002)
#    between the Butler et al. 2004 and Pakkala et al. 2002 references.
#
#    Building the new paragraph right after the Butler paragraph makes
#    this runtime inherit stray italic formatting (because that
#    paragraph contains an italic run somewhere inside it), so instead
#    we build the paragraph in a "clean" (non-italic) spot, apply the
#    correct italic run to the middle sentence only, and then move the
#    whole paragraph (cut/paste) to its correct location.
# ---------------------------------------------------------------------

# Paragraph "Förutom fridlysning ... (strikt skyddade djurarter)." has
# no italic runs at all, so inserting after it starts from a clean
# (non-italic) formatting state.
$cleanAnchorText = "Förutom fridlysning enligt §4 Artskyddsförordningen är tretåig hackspett även förtecknad i EU:s fågeldirektiv bilaga 1. Den ingår också i Natura 2000 och är förtecknad i Bernkonventionen bilaga II (strikt skyddade djurarter)."
$anchorPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.TrimEnd() -eq $cleanAnchorText) {
        $anchorPara = $d.Paragraphs($i)
        break
    }
}

$anchorPara.Range.InsertParagraphAfter()
$newParaIndex = $anchorPara.Index + 1
$newPara = $d.Paragraphs($newParaIndex)

$fullRefText = "Imbeau, L. & Desrochers, A. 2002. Foraging Ecology and Use of Drumming Trees by Three-Toed Woodpeckers. The Journal of Wildlife Management. Vol. 66, No. 1 (Jan., 2002), pp. 222-231."
$newPara.Range.Text = $fullRefText

$italicPart = "Foraging Ecology and Use of Drumming Trees by Three-Toed Woodpeckers. "
$newPara = $d.Paragraphs($newParaIndex)
$italicRange = $newPara.Range.Duplicate
$italicRange.Find.Execute($italicPart, $true, $false, $false, $false, $false, $true, 1, $false, "", 0) | Out-Null
$italicRange.Font.Italic = $true

# Cut the freshly built paragraph (with its ending paragraph mark) ...
$newPara = $d.Paragraphs($newParaIndex)
$newPara.Range.Cut() | Out-Null

# ... and find the Butler et al. reference paragraph to paste after it.
$butlerStart = "Butler, R., Angelstam, P., Ekelund, P. & Schlaeffer, R. 2004."
$butlerPara = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs($i).Range.Text.StartsWith($butlerStart)) {
        $butlerPara = $d.Paragraphs($i)
        break
    }
}

$pasteRange = $d.Range($butlerPara.Range.End, $butlerPara.Range.End)
$pasteRange.Paste() | Out-Null

# ---------------------------------------------------------------------
# 3) Update the date shown in the "first page" header.
# ---------------------------------------------------------------------
$headerRange = $d.Sections(1).Headers(2).Range
$headerRange.Find.Execute("2023-10-13", $true, $false, $false, $false, $false, $true, 1, $false, "2023-10-22", 2) | Out-Null
